# Replace the complex field (fldChar begin/.../instrText/fldChar end) that
# builds the "m:'doc.html'.fromHTMLURI()" expression with plain literal text
# runs wrapped in "{" / "}" - i.e. the field becomes inert template text
# instead of a live Word field, per TokenIteratorFieldRewriterSplit.

$d = $word.ActiveDocument

# Locate the paragraph that hosts the field (don't hard-code the index).
$targetParagraph = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $targetParagraph = $candidate
        break
    }
}

$range = $targetParagraph.Range

# Rebuild the paragraph's run content as literal text, keeping the existing
# "_GoBack" bookmark in place between "doc.html" and "'.fromHTMLURI()".
$openXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>{</w:t></w:r>
            <w:r><w:t>m</w:t></w:r>
            <w:r><w:t>:</w:t></w:r>
            <w:r><w:t>'</w:t></w:r>
            <w:r><w:t>doc.html</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t>'.fromHTMLURI()</w:t></w:r>
            <w:r><w:t xml:space="preserve">}</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$range.InsertXML($openXml) | Out-Null
